$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.298.46'
$ws.Range('E2').Value = '  +1.57%  '
$ws.Range('D3').Value = '3.903.41'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '526.60'
$ws.Range('E5').Value = '  +9.05%  '
$ws.Range('D6').Value = '144.81'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.614'
$ws.Range('E7').Value = '  -1.34%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -3.19%  '
$ws.Range('E10').Value = '  -5.03%  '
$ws.Range('E11').Value = '  -4.39%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '42.05'
$ws.Range('E12').Value = '  -2.46%  '
$ws.Range('D13').Value = '4.533.66'
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '10.28'
$ws.Range('E14').Value = '  -2.26%  '
$ws.Range('D15').Value = '3.926.30'
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('D16').Value = '1.23'
$ws.Range('E16').Value = '  +8.97%  '
$ws.Range('D17').Value = '13.99'
$ws.Range('E17').Value = '  -1.26%  '
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '19.74'
$ws.Range('E19').Value = '  -0.96%  '
$ws.Range('D20').Value = '69.256.59'
$ws.Range('E20').Value = '  +1.45%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '426.86'
$ws.Range('E21').Value = '  -0.61%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.37'
$ws.Range('E22').Value = '  -6.22%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '88.41'
$ws.Range('E23').Value = '  -0.78%  '
$ws.Range('D24').Value = '14.15'
$ws.Range('E24').Value = '  -4.40%  '
$ws.Range('E25').Value = '  +10.40%  '
$ws.Range('D26').Value = '11.43'
$ws.Range('E26').Value = '  -7.02%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.60'
$ws.Range('E27').Value = '  -3.83%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '36.45'
$ws.Range('E28').Value = '  -2.24%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '678.70'
$ws.Range('E29').Value = '  -4.66%  '
$ws.Range('D30').Value = '13.15'
$ws.Range('E30').Value = '  -2.35%  '
$ws.Range('E31').Value = '  -3.03%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.82'
$ws.Range('E32').Value = '  -2.89%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '68.58'
$ws.Range('E33').Value = '  +11.22%  '
$ws.Range('D34').Value = '0.0₃0882'
$ws.Range('E34').Value = '  +0.72%  '
$ws.Range('D35').Value = '0.435'
$ws.Range('E35').Value = '  +9.42%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.98'
$ws.Range('E36').Value = '  -1.12%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '39.97'
$ws.Range('E37').Value = '  -2.33%  '
$ws.Range('E38').Value = '  +1.88%  '
$ws.Range('E39').Value = '  +0.08%  '
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').Value = '3.25'
$ws.Range('E41').Value = '  +6.01%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0480'
$ws.Range('E42').Value = '  -3.53%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.17'
$ws.Range('E43').Value = '  +7.59%  '
$ws.Range('D44').Value = '2.79'
$ws.Range('E44').Value = '  -8.19%  '
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.35'
$ws.Range('E45').Value = '  -0.35%  '
$ws.Range('B46').Value = 'FLOKI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.000286'
$ws.Range('E46').Value = '  +18.93%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.140'
$ws.Range('E47').Value = '  -1.53%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.99'
$ws.Range('E48').Value = '  +7.01%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0349'
$ws.Range('E49').Value = '  -2.49%  '
$ws.Range('D50').Value = '2.743.69'
$ws.Range('E50').Value = '  +14.04%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '145.14'
$ws.Range('E51').Value = '  +0.25%  '
